$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44342
$ws.Range("M2").Value = 300

# Row 3
$ws.Range("D3").Value = 44319
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 120
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("R3").Value = "Provincia de Limarí"
$ws.Range("S3").Value = 1111
$ws.Range("T3").Value = 18

# Row 4
$ws.Range("D4").Value = 44355
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("S4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44358
$ws.Range("L5").Value = "Especial"
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18000
$ws.Range("P5").Value = 18000
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Provincia de Limarí"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 18

# Row 6
$ws.Range("D6").Value = 44358
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 100
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("S6").Value = 944

# Row 7
$ws.Range("D7").Value = 44294
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("Q7").Value = "$/caja 15 kilos granel"
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 800
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44326
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 300
$ws.Range("N8").Value = 20000
$ws.Range("O8").Value = 20000
$ws.Range("P8").Value = 20000
$ws.Range("S8").Value = 1111

# Row 9
$ws.Range("D9").Value = 44354
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 100

# Row 10
$ws.Range("D10").Value = 44348
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 200

# Row 13
$ws.Range("D13").Value = 44291
$ws.Range("L13").Value = "Primera"
$ws.Range("N13").Value = 12000
$ws.Range("O13").Value = 12000
$ws.Range("P13").Value = 12000
$ws.Range("Q13").Value = "$/caja 15 kilos granel"
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 800
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 44328
$ws.Range("L14").Value = "Especial"
$ws.Range("M14").Value = 250
$ws.Range("N14").Value = 20000
$ws.Range("O14").Value = 20000
$ws.Range("P14").Value = 20000
$ws.Range("S14").Value = 1111

# Row 15
$ws.Range("D15").Value = 44340
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 230
